$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply scraped cryptocurrency price/volume updates cell by cell.
# Numeric-looking text values are forced to remain text (matching the
# source data which stores prices/percentages as strings), then the
# cell style is restored to Normal so no stray formatting is left behind.

$ws.Range("D2").Value = '22.418.37'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '1.567.19'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3686'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.91'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3393'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.143'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07524'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.12'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.003'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.982'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.35%  '
$ws.Range("D16").Value = '1.572.03'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001119'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06792'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.355'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.25%  '
$ws.Range("D24").Value = '22.411.26'
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.362'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.647'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.049'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").Value = '1.747.43'
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.062'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.215'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.017'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.806'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08370'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02474'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.343'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2289'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06497'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.399'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6203'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.68%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.784'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5847'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.057'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '126.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.233'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07293'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.12%  '
